# feat: add 2022-Q3 data
#
# Before:  总计 | 2022-Q2 | 2022-Q1
# After:   总计 | 2022-Q3 | 2022-Q2 | 2022-Q1
#
# The existing "2022-Q2" sheet becomes the new "2022-Q3" sheet (refreshed with
# the latest quarter numbers), and a duplicate of its old data is inserted
# right after it, re-labelled "2022-Q2" (i.e. what used to be the Q2 tab now
# also exists, untouched, as a plain copy) followed by the unchanged "2022-Q1"
# tab.

$wb = $excel.ActiveWorkbook

$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# 1) Duplicate the current "2022-Q2" sheet, placing the copy right after the
#    original (i.e. before "2022-Q1"). Excel names it "2022-Q2 (2)".
$wsQ2.Copy($null, $wsQ2)
$wsQ2Copy = $wb.Worksheets.Item("2022-Q2 (2)")

# 2) Free up the "2022-Q2" name by promoting the original sheet to "2022-Q3",
#    then rename the duplicate (which still holds the old Q2 figures) back to
#    "2022-Q2".
$wsQ2.Name = "2022-Q3"
$wsQ2Copy.Name = "2022-Q2"

# 3) Refresh the figures on the (now) "2022-Q3" sheet with the new quarter's
#    data. D2:G2 are text-formatted figures (not numbers) in the source
#    sheet, so they're entered with a leading apostrophe to force text and
#    the cell style is then reset to "Normal" so no stray quote-prefix
#    formatting is left behind (matches the plain, unstyled source cells).
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Range("C2").Value = "嘉实全球房地产（QDII）"
$wsQ3.Range("D2").Value = "'0.38"
$wsQ3.Range("D2").Style = "Normal"
$wsQ3.Range("E2").Value = "'94.39"
$wsQ3.Range("E2").Style = "Normal"
$wsQ3.Range("F2").Value = "'2.32"
$wsQ3.Range("F2").Style = "Normal"
$wsQ3.Range("G2").Value = "'0.0088"
$wsQ3.Range("G2").Style = "Normal"
$wsQ3.Range("H2").Value = 10

# 4) Update the "总计" (totals) sheet: the old row for 2022-Q2 becomes the
#    2022-Q3 row, the old row for 2022-Q1 becomes the 2022-Q2 row (same
#    holding value as before, 0.01), and a new row is appended for 2022-Q1
#    with the value that used to belong to the old row 3 (0.02).
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("D3").Value = 0.01

# Copy row 3's formatting down to the new row 4 (matches A2/A3's bordered,
# bold-header style) before filling in row 4's own values.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.02

# Copying a sheet makes the new copy the active tab; restore "2022-Q1" (the
# tab that was active before this edit) as the selected sheet.
$wb.Worksheets.Item("2022-Q1").Activate()
